$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.435.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.572.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3730'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.89'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3392'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07579'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.138'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.997'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.967'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.581.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001120'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.88'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06739'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.286'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.439.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.335'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.644'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.10'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.007'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.758.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.048'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.167'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.971'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.787'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08385'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.378'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02473'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2283'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06517'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.457'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6226'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.94'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.809'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5794'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.074'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.216'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07315'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.24%  '
